# Atualização automática de preços de eletricidade
# Updates row 2 of the Spot_PT sheet with the newest daily hourly spot prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (date serial number)
$ws.Range("A2").Value = 45988

# Hourly prices 0h-1h ... 23h-24h (columns B..Y)
$ws.Range("B2").Value  = 88.89
$ws.Range("C2").Value  = 81.45999999999999
$ws.Range("D2").Value  = 73.81
$ws.Range("E2").Value  = 68.41
$ws.Range("F2").Value  = 66.38
$ws.Range("G2").Value  = 68.97
$ws.Range("H2").Value  = 77.45999999999999
$ws.Range("I2").Value  = 96.70999999999999
$ws.Range("J2").Value  = 99.89
$ws.Range("K2").Value  = 84.17
$ws.Range("L2").Value  = 55.37
$ws.Range("M2").Value  = 45.31
$ws.Range("N2").Value  = 45.91
$ws.Range("O2").Value  = 51.59
$ws.Range("P2").Value  = 55.48
$ws.Range("Q2").Value  = 63.83
$ws.Range("R2").Value  = 84.09999999999999
$ws.Range("S2").Value  = 101.56
$ws.Range("T2").Value  = 114.72
$ws.Range("U2").Value  = 119.01
$ws.Range("V2").Value  = 147.4
$ws.Range("W2").Value  = 155.93
$ws.Range("X2").Value  = 115.75
$ws.Range("Y2").Value  = 102.3

# Price_Daily_Avg
$ws.Range("Z2").Value  = 86.02

# Slot_4h_max / Slot_4h_price
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 130.34

# Slot_2h_frist / Slot_2h_frist_price
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 151.67

# Slot_2h_second / Slot_2h_second_price
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 116.86

# Slot_min_price (AG2) is unchanged by this update.
